$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update values per diff
$ws.Range("I7").Value = 1

$ws.Range("I8").Value = 1
$ws.Range("L8").Value = 1

$ws.Range("H9").Value = 30
$ws.Range("I9").Value = 1

$ws.Range("H10").Value = 31
$ws.Range("I10").Value = 1
$ws.Range("L10").Value = 1

$ws.Range("H11").Value = 32
$ws.Range("I11").Value = 1
$ws.Range("L11").Value = 0.1

# Update the active selection to L12 (was L13)
$ws.Range("L12").Select()
